# Applies the "Big Catch" copy refresh described in the commit diff.
#
# Two wrinkles forced the approach below:
#
# 1. Passing replacement text through Find.Execute's ReplaceWith arg (or
#    assigning Range.Text across a whole paragraph) runs it through Word's
#    smart-quote AutoCorrect, turning the straight apostrophe in
#    "Novomatic's" into a curly one. Locating the range with Find and then
#    writing through Range.InsertXML bypasses that.
#
# 2. Several target paragraphs start with an empty placeholder run
#    (<w:r/>) ahead of the text run. A plain Range.Text= (or a
#    Find.Execute replace) rewrites/collapses the whole paragraph and
#    silently drops that empty run. Rebuilding the exact paragraph XML
#    (pPr/rPr included) via Range.InsertXML keeps the run layout intact.
#    The very last paragraph in the body is a special case: handing
#    InsertXML that paragraph's full range (including its trailing
#    paragraph mark) duplicates the paragraph, so that mark is excluded
#    there and the already-present leading empty run is left alone
#    (re-adding it in the payload would duplicate it).

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Get-ParagraphByText($needle) {
    foreach ($para in $d.Paragraphs) {
        if ($para.Range.Text -like "*$needle*") {
            return $para
        }
    }
    return $null
}

function Wrap-Package($innerParagraphXml) {
    return "<?xml version='1.0' encoding='UTF-8' standalone='yes'?>" +
           "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
           "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
           "<pkg:xmlData><w:document xmlns:w='$wNs'><w:body>" +
           "<w:p>$innerParagraphXml</w:p>" +
           "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
}

function Set-ParagraphXml($needle, $innerParagraphXml) {
    $para = Get-ParagraphByText $needle
    if ($para -eq $null) {
        throw "Paragraph containing '$needle' not found"
    }

    $isLastParagraph = ($para.Range.End -ge $d.Content.End)
    if ($isLastParagraph) {
        # Exclude the trailing paragraph-mark character so InsertXML does
        # not duplicate the paragraph; the leading empty run already in
        # the document survives the edit on its own at this boundary, so
        # leave it out of the payload to avoid inserting a second copy.
        $target = $d.Range($para.Range.Start, $para.Range.End - 1)
        $target.InsertXML((Wrap-Package $innerParagraphXml))
    } else {
        $para.Range.InsertXML((Wrap-Package $innerParagraphXml))
    }
}

# --- Title / H1 heading (first occurrence) -------------------------------
# No sibling empty runs here, so a direct Find + Range.Text is safe and
# keeps things simple.
$rng = $d.Content
$found = $rng.Find.Execute("Play Big Catch Slot for Free - Game Review", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "Play Big Catch Free | Review of Novomatic's Slot Game"
}

# --- "What we like" bullets ------------------------------------------------
Set-ParagraphXml "Access to free spins and bonus games" ('<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Possibility of accessing free spins with a Bonus Game</w:t></w:r>')

Set-ParagraphXml "Structurally simple and minimalist design" ('<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Colorful and engaging underwater theme</w:t></w:r>')

Set-ParagraphXml "Cute and stylized symbols" ('<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Reliable and stable software for a smooth gaming experience</w:t></w:r>')

# --- "What we don't like" bullets ------------------------------------------
Set-ParagraphXml "Lacks special features and Autoplay" ('<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Lack of special features and Autoplay</w:t></w:r>')

Set-ParagraphXml "Theme may not appeal to all players" ('<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Limited bonus game with repetitive fishing hook mechanic</w:t></w:r>')

# --- Bold title repeated near the end of the document ----------------------
Set-ParagraphXml "Play Big Catch Slot for Free - Game Review" ('<w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Big Catch Free | Review of Novomatic' + [char]39 + 's Slot Game</w:t></w:r>')

# --- Italic meta-description paragraph (last paragraph in the body) --------
Set-ParagraphXml "Read our review of Big Catch, a simple yet engaging slot game by Novomatic with free spins and bonus features." ('<w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Big Catch, the underwater-themed slot game by Novomatic. Play it for free and enjoy simple gameplay mechanics.</w:t></w:r>')
